$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : game definition -------------------------------------------------
$ws.Range("A2").Value = "Math"
$ws.Range("B2").Value = "Add"
$ws.Range("C2").Value = "'"
$ws.Range("D2").Value = "'2"
$ws.Range("E2").Value = "'"
$ws.Range("E2").Value = ""

# --- Row 3 : MCQ question ------------------------------------------------------
$ws.Range("A3").Value = "Math"
$ws.Range("B3").Value = "Add2Numbers"
$ws.Range("C3").Value = "Add"
$ws.Range("D3").Value = "'1"
$ws.Range("E3").Value = "'MCQ"
$ws.Range("F3").Value = "12+1=?"
$ws.Range("G3").Value = "'10"
$ws.Range("H3").Value = "'1"
$ws.Range("I3").Value = "'12"
$ws.Range("J3").Value = "'13"
$ws.Range("K3").Value = "'13"
$ws.Range("L3").Value = "'10"

# --- Row 4 : True/False question ----------------------------------------------
$ws.Range("A4").Value = "Math"
$ws.Range("B4").Value = "Add"
$ws.Range("C4").Value = "Add2Numbers"
$ws.Range("D4").Value = "'2"
$ws.Range("E4").Value = "'T/F"
$ws.Range("F4").Value = "5+11=16"
$ws.Range("G4").Value = "'True"
$ws.Range("H4").Value = "'10"
$ws.Range("I4").Value = "'2+19=20"
$ws.Range("J4").Value = "'False"
$ws.Range("K4").Value = "'10"

# --- Final selection matches the authored workbook -----------------------------
$ws.Range("L3").Select()
